# Release v0.1.0-beta: Fix validation errors and update canonical URL
#
# Applies the metadata/content changes described by the commit:
#  - Version   1.0.0 -> 0.1.0
#  - Status    active -> draft
#  - Date      2025-11-28T01:24:36+00:00 -> 2025-12-26T14:13:58+00:00
#  - Description (was blank) -> filled in with the extension's description
#  - Elements table: the root "Extension" element's Definition cell is
#    updated from the placeholder "An Extension" text to the same new
#    description text used on the Metadata sheet.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

$description = "Extension to link nursing interventions to the patient goals they are intended to achieve. Supports goal-directed care planning and intervention tracking."

# --- Metadata sheet (Property / Value pairs) ---
$metadata.Range("B3").Value = "0.1.0"
$metadata.Range("B6").Value = "draft"
$metadata.Range("B8").Value = "2025-12-26T14:13:58+00:00"
$metadata.Range("B11").Value = $description

# --- Elements sheet (row 2 = root "Extension" element, column M = Definition) ---
$elements.Range("M2").Value = $description
